{"js": "// Split the single \"#7cc867#fb5b89#c885da#f9cd59\" paragraph (under the\n// \"Highlights\" heading) into four paragraphs, one per highlight color,\n// each annotated with its count: \"#COLOR: N\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = \"#7cc867#fb5b89#c885da#f9cd59\";\nconst replacements = [\"#7cc867: 31\", \"#fb5b89: 50\", \"#c885da: 19\", \"#f9cd59: 33\"];\n\nconst hit = paragraphs.items.find((p) => p.text === target);\nif (hit) {\n  // Rewrite the original paragraph with the first highlight's line\u2026\n  hit.insertText(replacements[0], \"Replace\");\n  // \u2026then insert the remaining lines as new paragraphs right after it,\n  // each chained off the previous so they land in order.\n  let anchor = hit;\n  for (let i = 1; i < replacements.length; i++) {\n    anchor = anchor.insertParagraph(replacements[i], \"After\");\n  }\n  await context.sync();\n}\n", "ps1": "# Split the single \"#7cc867#fb5b89#c885da#f9cd59\" paragraph (under the\n# \"Highlights\" heading) into four paragraphs, one per highlight color,\n# each annotated with its count: \"#COLOR: N\".\n$d = $word.ActiveDocument\n\n$target = \"#7cc867#fb5b89#c885da#f9cd59\"\n$cr = [char]13\n$replacement = \"#7cc867: 31\" + $cr + \"#fb5b89: 50\" + $cr + \"#c885da: 19\" + $cr + \"#f9cd59: 33\"\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    if ($r.Text -eq ($target + $cr)) {\n        # Exclude the trailing paragraph mark from the replaced span, then\n        # assign text with embedded carriage returns so Word splits it into\n        # four separate paragraphs in place of the original one.\n        $body = $r.Duplicate\n        [void]$body.MoveEnd(1, -1)\n        $body.Text = $replacement\n        break\n    }\n}\n"}
